# This script applies a row-rotation edit to the "Artfynd" sheet.
# The underlying change (per the source diff) moves data so that:
#   - the record that was in row 4 becomes row 2
#   - the record that was in row 2 becomes row 3
#   - the record that was in row 3 becomes row 4
# i.e. a cyclic shift of the three data rows (2,3,4).
# Below we set each cell to its target value directly (and clear cells that
# should become empty), which reproduces that row rotation.
# For cells whose new text looks like a number or an ISO date (e.g. "1" or
# "2022-01-26") we force the cell to Text format first, so Excel keeps the
# value as literal text instead of auto-converting it to a number/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: set new values (rotated from original rows)
$ws.Range("A2").Value = 98350752
$ws.Range("B2").Value = 93132
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 2671
$ws.Range("F2").Value = "Fällmossa"
$ws.Range("G2").Value = "Antitrichia curtipendula"
$ws.Range("H2").Value = "(Hedw.) Brid."
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("P2").Value = "Korpklint, Ög"
$ws.Range("Q2").Value = 573934.5192830344
$ws.Range("R2").Value = 6505660.931744166
$ws.Range("S2").Value = 15
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2022-01-26"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2022-01-26"
$ws.Range("AC2").ClearContents()
$ws.Range("AI2").ClearContents()
$ws.Range("AW2").Value = "Marika Sjödin"
$ws.Range("AX2").Value = "Marika Sjödin, Eva Siljeholm"
$ws.Range("AY2").ClearContents()

# Row 3: set new values (rotated from original rows)
$ws.Range("A3").Value = 80000956
$ws.Range("B3").Value = 57133
$ws.Range("D3").Value = "VU"
$ws.Range("E3").Value = 100041
$ws.Range("F3").Value = "Hasselsnok"
$ws.Range("G3").Value = "Coronella austriaca"
$ws.Range("H3").Value = "Laurenti, 1768"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "ex."
$ws.Range("N3").Value = "observerad"
$ws.Range("P3").Value = "Ekbacken, Ög"
$ws.Range("Q3").Value = 574245.9331973131
$ws.Range("R3").Value = 6505393.860064601
$ws.Range("S3").Value = 5
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2015-07-13"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2015-07-13"
$ws.Range("AC3").Value = "observerad"
$ws.Range("AI3").Value = "rasbrant"
$ws.Range("AW3").Value = "Elin Håkansson"
$ws.Range("AX3").Value = "Mikael Hagström"
$ws.Range("AY3").Value = "Ostlänken Norrköpings kommun (OLP2)"

# Row 4: set new values (rotated from original rows)
$ws.Range("A4").Value = 97784188
$ws.Range("B4").Value = 78098
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6453
$ws.Range("F4").Value = "Vedskivlav"
$ws.Range("G4").Value = "Hertelidea botryosa"
$ws.Range("H4").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q4").Value = 574267.5050753297
$ws.Range("R4").Value = 6505328.227546699
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2021-12-30"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2021-12-30"
